$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Neg_Change" (sheet1): replace the 10 old data rows (rows
# 2-11) with the 4 new data rows, then drop the now-unused rows.
# -----------------------------------------------------------------
$wsNeg = $wb.Worksheets.Item("Neg_Change")

$negData = @(
    @("GODREJPROP", 2324, 2330, 2275.3, 2295, 394259, 949508, -0.5847754837242024, "GODREJPROP"),
    @("OBEROIRLTY", 1737, 1740, 1703.1, 1706.2, 428250, 1003301, -0.5731590021339558, "OBEROIRLTY"),
    @("IDEA", 10.1, 10.15, 9.380000000000001, 9.44, 1456935645, 3237582345, -0.5499927137760569, "IDEA"),
    @("RBLBANK", 323, 326.35, 319.15, 324.95, 8631329, 18611108, -0.5362270209812333, "RBLBANK")
)

for ($i = 0; $i -lt $negData.Length; $i++) {
    $row = $negData[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $wsNeg.Cells.Item(2 + $i, 1 + $j).Value = $row[$j]
    }
}

# Remove the leftover rows (old rows 6-11) that have no replacement data.
$wsNeg.Range("A6:I11").EntireRow.Delete()

# -----------------------------------------------------------------
# Sheet "Pos_Change" (sheet2): replace the 16 old data rows (rows
# 2-17) with the 17 new data rows (rows 2-18).
# -----------------------------------------------------------------
$wsPos = $wb.Worksheets.Item("Pos_Change")

$posData = @(
    @("BAJAJ-AUTO", 9100, 9130, 9019, 9057, 459991, 310061, 0.4835500111268428, "BAJAJ-AUTO"),
    @("ITC", 420.5, 421.65, 415.9, 417.85, 17163147, 11603679, 0.4791125297416449, "ITC"),
    @("COALINDIA", 396.95, 398.2, 389.35, 391.5, 7672258, 5438773, 0.4106597204921036, "COALINDIA"),
    @("TRENT", 4800.6, 4828.4, 4692.5, 4733, 616890, 414172, 0.4894536569348, "TRENT"),
    @("CANBK", 130.39, 131.77, 127.75, 130.1, 44849817, 30460506, 0.4723923824509021, "CANBK"),
    @("TATAPOWER", 400.25, 401.85, 396.1, 399.05, 3949368, 2799280, 0.4108513617787431, "TATAPOWER"),
    @("CGPOWER", 730, 731.45, 717.4, 722, 2309934, 1573062, 0.4684316320653604, "CGPOWER"),
    @("HINDZINC", 481.4, 481.4, 468.2, 471.5, 5637647, 3767112, 0.4965435060067235, "HINDZINC"),
    @("UNIONBANK", 145.9, 147.73, 144.7, 147.4, 15194370, 10145527, 0.4976422614616274, "UNIONBANK"),
    @("NYKAA", 256, 259.5, 253.72, 258.77, 6379166, 4262941, 0.496423713112614, "NYKAA"),
    @("JUBLFOOD", 595.7, 601.65, 592.1, 601.2, 1169233, 804991, 0.4524795929395484, "JUBLFOOD"),
    @("APLAPOLLO", 1766.8, 1793.1, 1753.1, 1776, 438596, 280131, 0.5656817703146028, "APLAPOLLO"),
    @("IGL", 213.47, 213.89, 208.89, 210.58, 2413492, 1557159, 0.5499329227137369, "IGL"),
    @("IIFL", 506, 520, 505.85, 513.3, 2383613, 1658725, 0.4370151773199294, "IIFL"),
    @("TATACHEM", 901, 912.1, 901, 903.85, 453366, 314412, 0.4419487805808939, "TATACHEM"),
    @("MANAPPURAM", 275.8, 279.8, 274, 277.65, 3653581, 2290776, 0.5949097598368414, "MANAPPURAM"),
    @("CROMPTON", 292.8, 293, 287.5, 289.4, 4533574, 3112928, 0.4563696943841939, "CROMPTON")
)

for ($i = 0; $i -lt $posData.Length; $i++) {
    $row = $posData[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $wsPos.Cells.Item(2 + $i, 1 + $j).Value = $row[$j]
    }
}
